$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.072918772697449
$ws.Range("B1").Value = 1.712220191955566
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.866223335266113
$ws.Range("E1").Value = 1.157050013542175
